$wb = $excel.ActiveWorkbook

# New row data (row 77) for each of the 4 sheets, in workbook tab order:
#   1 FE_LFT_#1 , 2 FE_LFT_#2 , 3 FE_PLT_#1 , 4 FE_PLT_#2
# Columns: A time, B total-len(hex), C id(hex), D actual-len(hex), E checksum(hex),
#          F total-len dec, G id dec, H actual-len dec, I checksum dec

$timeVal = 45863.49141203704

# G values are written as the exact decimal expansion of the target IEEE-754
# double (computed offline) so the stored bit pattern matches exactly -
# straightforward decimal-exponent literals (e.g. 7.59e+23) aren't accepted
# by the script parser, and multiplying by a power of ten can round
# differently in the last bit.
$rows = @(
    @{ B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x24"; E = "0xf"; F = 380; G = 759863127514710945038336.0; H = 292; I = 15 },
    @{ B = "0x01,0x90"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x34"; E = "0xe"; F = 400; G = 568432987514711010443264.0; H = 308; I = 14 },
    @{ B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x61"; E = "0x3"; F = 110; G = 568631262647113970876416.0; H = 97;  I = 3 },
    @{ B = "0x00,0x6e"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x60"; E = "0x3"; F = 110; G = 985046333984776009023488.0; H = 96;  I = 3 }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rows[$i - 1]

    $newRow = 77

    # Match the date/time number format used by the rest of column A (style "2").
    $ws.Range("A$newRow").NumberFormat = $ws.Range("A" + ($newRow - 1)).NumberFormat
    $ws.Range("A$newRow").Value2 = $timeVal

    $ws.Range("B$newRow").Value2 = $data.B
    $ws.Range("C$newRow").Value2 = $data.C
    $ws.Range("D$newRow").Value2 = $data.D
    $ws.Range("E$newRow").Value2 = $data.E
    $ws.Range("F$newRow").Value2 = $data.F
    $ws.Range("G$newRow").Value2 = $data.G
    $ws.Range("H$newRow").Value2 = $data.H
    $ws.Range("I$newRow").Value2 = $data.I
}
